# Updated version of HSI, Review and SIQ, as well as the RTM itself.
#
# The review points on the "CYRS" sheet that previously had no
# Acceptance status (column D) are now marked "Accepted".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CYRS")

$ws.Range("D2").Value = "Accepted"
$ws.Range("D5").Value = "Accepted"
$ws.Range("D6").Value = "Accepted"
$ws.Range("D7").Value = "Accepted"
$ws.Range("D8").Value = "Accepted"
$ws.Range("D9").Value = "Accepted"

# Move the active selection to D2 (matches the saved selection state)
$ws.Activate() | Out-Null
$ws.Range("D2").Select() | Out-Null
